$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# Update product name in A2
$ws.Range("A2").Value = "Moto"

# Update quantity/value in D2
$ws.Range("D2").Value = 1

# Move the active selection to A2 (matches the saved view state)
$ws.Range("A2").Select()
